# hp and mp bar text
# Add two new rows (28 and 29) to the "Todo " sheet describing new tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Todo ")

# Row 28: new task "min damage 1"
$ws.Cells.Item(28, 1).Value = "min damage 1"   # A28 - Issue
$ws.Cells.Item(28, 4).Value = "Fish"           # D28 - Creator
$ws.Cells.Item(28, 5).Value = "Fish"           # E28 - PIC
$ws.Cells.Item(28, 6).Value = "14 Jan"         # F28 - Create date

# Row 29: new task "db and save"
$ws.Cells.Item(29, 1).Value = "db and save"    # A29 - Issue
$ws.Cells.Item(29, 4).Value = "Fish"           # D29 - Creator
$ws.Cells.Item(29, 5).Value = "Fish"           # E29 - PIC

# Match the column F style (text number format) used on other rows,
# so the new Create date value stays formatted consistently.
$ws.Range("F28").NumberFormat = $ws.Range("F27").NumberFormat

# Update the view to match the author's final selection/scroll position.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E28").Select()
